# "Add comparison to report EFs"
#
# 1. Text/label clean-up in the shared strings used by the Slurry and
#    Application sheets (capitalisation + spelling fixes, and clearer
#    English labels for the application-method codes).
# 2. Make "Application" the active/selected sheet (was "Slurry"), with a
#    fresh cell selection on each of the two affected sheets.
# 3. Widen column A on the "Application" sheet so the new, longer labels
#    are readable, and tighten the row height on the two rows whose
#    labels changed.

$wb = $excel.ActiveWorkbook

$slurry = $wb.Worksheets.Item("Slurry")
$climate = $wb.Worksheets.Item("Climate")
$app = $wb.Worksheets.Item("Application")

# --- 1. Shared-string text fixes -------------------------------------------------

# Slurry sheet: man.source column + the (duplicated) manure-type names.
$slurry.Range("B2").Value = "Pig"
$slurry.Range("B3").Value = "Cattle"
$slurry.Range("A4").Value = "Afgasset biomasse"
$slurry.Range("B4").Value = "Digestate"
$slurry.Range("B5").Value = "Pig"
$slurry.Range("B6").Value = "Cattle"
$slurry.Range("A7").Value = "Afgasset biomasse"
$slurry.Range("B7").Value = "Digestate"

# Application sheet: app.mthd / incorp columns.
$app.Range("A2").Value = "Trailing hose"
$app.Range("B2").Value = "None"
$app.Range("A3").Value = "Trailing hose"
$app.Range("B3").Value = "Shallow"
$app.Range("A4").Value = "Trailing hose"
$app.Range("B4").Value = "Deep"
$app.Range("A5").Value = "Open slot injection"
$app.Range("B5").Value = "None"
$app.Range("A6").Value = "Closed slot injection"
$app.Range("B6").Value = "None"

# --- 2. Active sheet / selection -------------------------------------------------

# "Application" becomes the active tab (was "Slurry"); give both sheets a
# fresh selection matching the new state of the workbook.
$slurry.Activate()
$slurry.Range("B8").Select() | Out-Null

$app.Activate()
$app.Range("B7").Select() | Out-Null

# --- 3. Column width / row height on the Application sheet ----------------------

# Column A needs to be wide enough for the new English labels.
$app.Columns.Item(1).ColumnWidth = 16.57142857142857

# Rows 3 & 4 (Shallow / Deep) tighten up to the sheet's normal row height.
$app.Rows.Item(3).RowHeight = 12.8
$app.Rows.Item(4).RowHeight = 12.8
